$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.397.56"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "3.046.92"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.206"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("D10").Value = "3.045.84"
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.438"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.24"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.82%  "
$ws.Range("D14").Value = "3.604.11"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.83"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "76.307.36"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "3.037.99"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  +3.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.34"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.76"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000108"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.28"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "495.31"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "20.60"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.15"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.118"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.03"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "192.28"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.379"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.797"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +20.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.09"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.90"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.14%  "
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.41"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.599"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.88"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.18%  "
